$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 244.36
$ws.Range("I19").Value = 206.16667
$ws.Range("K19").Value = 206.16667
$ws.Range("M19").Value = -31.16667000000001
$ws.Range("H96").Value = 544.6
$ws.Range("I96").Value = 548.5
$ws.Range("J96").Value = 529
$ws.Range("K96").Value = 1645.5
$ws.Range("L96").Value = 1587
$ws.Range("M96").Value = -272.5
$ws.Range("N96").Value = -4333
$ws.Range("H98").Value = 1081.5
$ws.Range("I98").Value = 1090.5555
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 1090.5555
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = 407.4445000000001
$ws.Range("N98").Value = -3996
$ws.Range("H116").Value = 13022.777
$ws.Range("I116").Value = 34668.332
$ws.Range("J116").Value = 2200
$ws.Range("K116").Value = 34668.332
$ws.Range("L116").Value = 2200
$ws.Range("M116").Value = -31226.332
$ws.Range("N116").Value = -9084
$ws.Range("H122").Value = 1081.5
$ws.Range("I122").Value = 1090.5555
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 3271.6665
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -821.6664999999998
$ws.Range("N122").Value = -7900

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9537.429
$ws.Range("I61").Value = 10749.917
$ws.Range("J61").Value = 2262.5
$ws.Range("K61").Value = 10749.917
$ws.Range("L61").Value = 2262.5
$ws.Range("M61").Value = -10537.917
$ws.Range("N61").Value = -2686.5
$ws.Range("H74").Value = 1500.0317
$ws.Range("I74").Value = 1439.1923
$ws.Range("J74").Value = 1787.6364
$ws.Range("K74").Value = 1439.1923
$ws.Range("L74").Value = 1787.6364
$ws.Range("M74").Value = -565.1922999999999
$ws.Range("N74").Value = -3535.6364
$ws.Range("H77").Value = 1500.0317
$ws.Range("I77").Value = 1439.1923
$ws.Range("J77").Value = 1787.6364
$ws.Range("K77").Value = 7195.961499999999
$ws.Range("L77").Value = 8938.182000000001
$ws.Range("M77").Value = -2827.961499999999
$ws.Range("N77").Value = -17674.182
$ws.Range("H97").Value = 866.86664
$ws.Range("J97").Value = 1388.8334
$ws.Range("L97").Value = 1388.8334
$ws.Range("N97").Value = -2380.8334
$ws.Range("H132").Value = 9939.474
$ws.Range("I132").Value = 2050.8572
$ws.Range("J132").Value = 14541.167
$ws.Range("K132").Value = 6152.571599999999
$ws.Range("L132").Value = 43623.501
$ws.Range("M132").Value = -3622.571599999999
$ws.Range("N132").Value = -48683.501
$ws.Range("H136").Value = 9537.429
$ws.Range("I136").Value = 10749.917
$ws.Range("J136").Value = 2262.5
$ws.Range("K136").Value = 32249.751
$ws.Range("L136").Value = 6787.5
$ws.Range("M136").Value = -29699.751
$ws.Range("N136").Value = -11887.5

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2214.4211
$ws.Range("I94").Value = 1767.375
$ws.Range("J94").Value = 2539.5454
$ws.Range("K94").Value = 1767.375
$ws.Range("L94").Value = 2539.5454
$ws.Range("M94").Value = -1316.375
$ws.Range("N94").Value = -3441.5454
$ws.Range("H107").Value = 1050.45
$ws.Range("I107").Value = 1090.1538
$ws.Range("J107").Value = 976.7143
$ws.Range("K107").Value = 1090.1538
$ws.Range("L107").Value = 976.7143
$ws.Range("M107").Value = 829.8462
$ws.Range("N107").Value = -4816.7143
$ws.Range("H134").Value = 5000.0835
$ws.Range("I134").Value = 5912.36
$ws.Range("K134").Value = 17737.08
$ws.Range("M134").Value = -15202.08

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5344.971
$ws.Range("I31").Value = 841.52
$ws.Range("J31").Value = 16603.6
$ws.Range("K31").Value = 841.52
$ws.Range("L31").Value = 16603.6
$ws.Range("M31").Value = -546.52
$ws.Range("N31").Value = -17193.6
$ws.Range("H34").Value = 5344.971
$ws.Range("I34").Value = 841.52
$ws.Range("J34").Value = 16603.6
$ws.Range("K34").Value = 841.52
$ws.Range("L34").Value = 16603.6
$ws.Range("M34").Value = -639.52
$ws.Range("N34").Value = -17007.6
$ws.Range("H58").Value = 1145.5
$ws.Range("I58").Value = 920.6087
$ws.Range("J58").Value = 2180
$ws.Range("K58").Value = 920.6087
$ws.Range("L58").Value = 2180
$ws.Range("M58").Value = -717.6087
$ws.Range("N58").Value = -2586
$ws.Range("H132").Value = 2215.0715
$ws.Range("I132").Value = 1421.4
$ws.Range("J132").Value = 4199.25
$ws.Range("K132").Value = 4264.200000000001
$ws.Range("L132").Value = 12597.75
$ws.Range("M132").Value = -1734.200000000001
$ws.Range("N132").Value = -17657.75
$ws.Range("H134").Value = 1912.2
$ws.Range("I134").Value = 1890.25
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 5670.75
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -3135.75
$ws.Range("N134").Value = -11070
$ws.Range("H136").Value = 1145.5
$ws.Range("I136").Value = 920.6087
$ws.Range("J136").Value = 2180
$ws.Range("K136").Value = 2761.8261
$ws.Range("L136").Value = 6540
$ws.Range("M136").Value = -211.8261000000002
$ws.Range("N136").Value = -11640

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 2962.6
$ws.Range("I76").Value = 1406.5
$ws.Range("K76").Value = 4219.5
$ws.Range("M76").Value = -3836.5
$ws.Range("H79").Value = 2962.6
$ws.Range("I79").Value = 1406.5
$ws.Range("K79").Value = 4219.5
$ws.Range("M79").Value = -2893.5
$ws.Range("H82").Value = 3683.3333
$ws.Range("I82").Value = 1150
$ws.Range("J82").Value = 4000
$ws.Range("K82").Value = 3450
$ws.Range("L82").Value = 12000
$ws.Range("M82").Value = -3044
$ws.Range("N82").Value = -12812
$ws.Range("H85").Value = 3683.3333
$ws.Range("I85").Value = 1150
$ws.Range("J85").Value = 4000
$ws.Range("K85").Value = 3450
$ws.Range("L85").Value = 12000
$ws.Range("M85").Value = -2046
$ws.Range("N85").Value = -14808
$ws.Range("H88").Value = 4286.364
$ws.Range("J88").Value = 4286.364
$ws.Range("L88").Value = 12859.092
$ws.Range("N88").Value = -13715.092
$ws.Range("H91").Value = 4286.364
$ws.Range("J91").Value = 4286.364
$ws.Range("L91").Value = 12859.092
$ws.Range("N91").Value = -15823.092
$ws.Range("H94").Value = 3476.923
$ws.Range("J94").Value = 3476.923
$ws.Range("L94").Value = 10430.769
$ws.Range("N94").Value = -11782.769
$ws.Range("H97").Value = 16667217
$ws.Range("H100").Value = 3480
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("H129").Value = 15152923
$ws.Range("I129").Value = 22223074
$ws.Range("J129").Value = 2599.8572
$ws.Range("K129").Value = 66669222
$ws.Range("L129").Value = 7799.571599999999
$ws.Range("M129").Value = -66664222
$ws.Range("N129").Value = -17799.5716
$ws.Range("H137").Value = 10753
$ws.Range("I137").Value = 16488.428
$ws.Range("J137").Value = 7885.2856
$ws.Range("K137").Value = 49465.284
$ws.Range("L137").Value = 23655.8568
$ws.Range("M137").Value = -44365.284
$ws.Range("N137").Value = -33855.8568

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6124.5806
$ws.Range("I70").Value = 6113.923
$ws.Range("J70").Value = 6180
$ws.Range("K70").Value = 6113.923
$ws.Range("L70").Value = 6180
$ws.Range("M70").Value = -5843.923
$ws.Range("N70").Value = -6720
$ws.Range("H73").Value = 6124.5806
$ws.Range("I73").Value = 6113.923
$ws.Range("J73").Value = 6180
$ws.Range("K73").Value = 6113.923
$ws.Range("L73").Value = 6180
$ws.Range("M73").Value = -5177.923
$ws.Range("N73").Value = -8052
$ws.Range("H80").Value = 3598.75
$ws.Range("I80").Value = 3997.5
$ws.Range("K80").Value = 3997.5
$ws.Range("M80").Value = -2999.5
$ws.Range("H83").Value = 3598.75
$ws.Range("I83").Value = 3997.5
$ws.Range("K83").Value = 19987.5
$ws.Range("M83").Value = -14995.5
$ws.Range("H97").Value = 673.7646999999999
$ws.Range("I97").Value = 702.0909
$ws.Range("J97").Value = 621.8333
$ws.Range("K97").Value = 702.0909
$ws.Range("L97").Value = 621.8333
$ws.Range("M97").Value = -206.0909
$ws.Range("N97").Value = -1613.8333
$ws.Range("H113").Value = 333334080
$ws.Range("I113").Value = 1000000000
$ws.Range("J113").Value = 1100
$ws.Range("K113").Value = 1000000000
$ws.Range("L113").Value = 1100
$ws.Range("M113").Value = -999997830
$ws.Range("N113").Value = -5440
$ws.Range("H132").Value = 3647.8647
$ws.Range("I132").Value = 3793.8462
$ws.Range("J132").Value = 3568.7917
$ws.Range("K132").Value = 11381.5386
$ws.Range("L132").Value = 10706.3751
$ws.Range("M132").Value = -8851.5386
$ws.Range("N132").Value = -15766.3751

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2181.6086
$ws.Range("I61").Value = 1885.1333
$ws.Range("J61").Value = 2737.5
$ws.Range("K61").Value = 1885.1333
$ws.Range("L61").Value = 2737.5
$ws.Range("M61").Value = -1683.1333
$ws.Range("N61").Value = -3141.5
$ws.Range("H93").Value = 33348000
$ws.Range("I93").Value = 26062.5
$ws.Range("J93").Value = 71430216
$ws.Range("K93").Value = 26062.5
$ws.Range("L93").Value = 71430216
$ws.Range("M93").Value = -24814.5
$ws.Range("N93").Value = -71432712
$ws.Range("H113").Value = 2181.6086
$ws.Range("I113").Value = 1885.1333
$ws.Range("J113").Value = 2737.5
$ws.Range("K113").Value = 1885.1333
$ws.Range("L113").Value = 2737.5
$ws.Range("M113").Value = 284.8667
$ws.Range("N113").Value = -7077.5
$ws.Range("H122").Value = 5104719.5
$ws.Range("I122").Value = 6496098
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 19488294
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -19485844
$ws.Range("N122").Value = -13900

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 42750
$ws.Range("J98").Value = 42750
$ws.Range("L98").Value = 42750
$ws.Range("N98").Value = -48740
$ws.Range("H107").Value = 55555972
$ws.Range("I107").Value = 83333740
$ws.Range("J107").Value = 423.16666
$ws.Range("K107").Value = 250001220
$ws.Range("L107").Value = 1269.49998
$ws.Range("M107").Value = -249999300
$ws.Range("N107").Value = -5109.499980000001
$ws.Range("H136").Value = 1153.0769
$ws.Range("I136").Value = 769
$ws.Range("J136").Value = 2433.3333
$ws.Range("K136").Value = 2307
$ws.Range("L136").Value = 7299.999899999999
$ws.Range("M136").Value = 243
$ws.Range("N136").Value = -12399.9999

# Remove M100 cell entirely (CUL sheet)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M100").ClearContents()
